# Add a new "Save" column (H) to the s_vals sheet, mirroring the
# header style used by the existing columns (B1:G1) and filling the
# data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold font, border, centered
# alignment) from the last existing header cell (G1) onto the new
# header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column's data rows with 0.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
